$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..141 (140 rows) get new EPC (col B) and Serial (col D) values.
# EPC is a 24-char hex string starting at 30300E890A01810077359401, incrementing
# by 1 per row. Serial starts at 2000000001, incrementing by 1 per row.

$epcPrefix = "30300E890A018100773594"
$startSuffix = 0x01
$startSerial = 1

for ($i = 0; $i -lt 140; $i++) {
    $row = $i + 2
    $suffixVal = $startSuffix + $i
    $hexSuffix = "{0:X2}" -f $suffixVal
    $epc = $epcPrefix + $hexSuffix
    $serial = "{0:D10}" -f (2000000000 + $startSerial + $i)

    $bCell = $ws.Cells.Item($row, 2)
    $dCell = $ws.Cells.Item($row, 4)

    # Serial is a purely-numeric-looking string; force text formatting so it
    # keeps its original string type instead of being coerced to a number.
    $dCell.NumberFormat = "@"

    $bCell.Value = $epc
    $dCell.Value = $serial
}
